$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dtype")

# Insert a new row at row 71 (shifts rows 71.. down by one).
$ws.Rows("71:71").Insert()

# Copy the formatting of the row above (row 70), which already matches the
# formatting used throughout this "OSeMOSYS-UGA inputs" section, into the
# freshly inserted row.
$ws.Range("A70:C70").Copy()
$ws.Range("A71:C71").PasteSpecial(-4122)

# Fill in the new parameter entry: AccumulatedAnnualDemand (float).
$ws.Range("A71").Value = "OSeMOSYS-UGA inputs"
$ws.Range("B71").Value = "AccumulatedAnnualDemand"
$ws.Range("C71").Value = "float"

$ws.Activate()
